$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 466
$ws.Range("J49").Value = 499
$ws.Range("L49").Value = 1497
$ws.Range("N49").Value = -1769
$ws.Range("H64").Value = 4316.6665
$ws.Range("I64").Value = 4990
$ws.Range("J64").Value = 3980
$ws.Range("K64").Value = 4990
$ws.Range("L64").Value = 3980
$ws.Range("M64").Value = -4742
$ws.Range("N64").Value = -4476
$ws.Range("H67").Value = 4316.6665
$ws.Range("I67").Value = 4990
$ws.Range("J67").Value = 3980
$ws.Range("K67").Value = 4990
$ws.Range("L67").Value = 3980
$ws.Range("M67").Value = -4132
$ws.Range("N67").Value = -5696
$ws.Range("H74").Value = 2666.6667
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 2000
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -3064
$ws.Range("N74").Value = -3872
$ws.Range("H76").Value = 4998.5454
$ws.Range("I76").Value = 5220
$ws.Range("J76").Value = 4002
$ws.Range("K76").Value = 5220
$ws.Range("L76").Value = 4002
$ws.Range("M76").Value = -4905
$ws.Range("N76").Value = -4632
$ws.Range("H77").Value = 2666.6667
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -15320
$ws.Range("N77").Value = -19360
$ws.Range("H79").Value = 4998.5454
$ws.Range("I79").Value = 5220
$ws.Range("J79").Value = 4002
$ws.Range("K79").Value = 5220
$ws.Range("L79").Value = 4002
$ws.Range("M79").Value = -4128
$ws.Range("N79").Value = -6186
$ws.Range("H132").Value = 9266626
$ws.Range("I132").Value = 13340113
$ws.Range("J132").Value = 8702.091
$ws.Range("K132").Value = 40020339
$ws.Range("L132").Value = 26106.273
$ws.Range("M132").Value = -40017809
$ws.Range("N132").Value = -31166.273
$ws.Range("H137").Value = 1518.6774
$ws.Range("I137").Value = 1142.5294
$ws.Range("J137").Value = 1975.4286
$ws.Range("K137").Value = 3427.5882
$ws.Range("L137").Value = 5926.2858
$ws.Range("M137").Value = -877.5881999999997
$ws.Range("N137").Value = -11026.2858
$ws.Range("H138").Value = 1691.909
$ws.Range("I138").Value = 1659.2
$ws.Range("J138").Value = 1696.1025
$ws.Range("K138").Value = 4977.6
$ws.Range("L138").Value = 5088.3075
$ws.Range("M138").Value = 162.3999999999996
$ws.Range("N138").Value = -15368.3075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1598.05
$ws.Range("I32").Value = 1426.3586
$ws.Range("J32").Value = 3572.5
$ws.Range("K32").Value = 1426.3586
$ws.Range("L32").Value = 3572.5
$ws.Range("M32").Value = -1139.3586
$ws.Range("N32").Value = -4146.5
$ws.Range("H74").Value = 1578.75
$ws.Range("I74").Value = 1166.1538
$ws.Range("J74").Value = 3366.6667
$ws.Range("K74").Value = 1166.1538
$ws.Range("L74").Value = 3366.6667
$ws.Range("M74").Value = -292.1538
$ws.Range("N74").Value = -5114.6667
$ws.Range("H77").Value = 1578.75
$ws.Range("I77").Value = 1166.1538
$ws.Range("J77").Value = 3366.6667
$ws.Range("K77").Value = 5830.769
$ws.Range("L77").Value = 16833.3335
$ws.Range("M77").Value = -1462.769
$ws.Range("N77").Value = -25569.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2819.7144
$ws.Range("I20").Value = 2550.875
$ws.Range("J20").Value = 3680
$ws.Range("K20").Value = 2550.875
$ws.Range("L20").Value = 3680
$ws.Range("M20").Value = -2303.875
$ws.Range("N20").Value = -4174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1503.4667
$ws.Range("I31").Value = 1390.1842
$ws.Range("J31").Value = 2118.4285
$ws.Range("K31").Value = 1390.1842
$ws.Range("L31").Value = 2118.4285
$ws.Range("M31").Value = -1095.1842
$ws.Range("N31").Value = -2708.4285
$ws.Range("H34").Value = 1503.4667
$ws.Range("I34").Value = 1390.1842
$ws.Range("J34").Value = 2118.4285
$ws.Range("K34").Value = 1390.1842
$ws.Range("L34").Value = 2118.4285
$ws.Range("M34").Value = -1188.1842
$ws.Range("N34").Value = -2522.4285
$ws.Range("H62").Value = 6898791
$ws.Range("I62").Value = 2329.5454
$ws.Range("J62").Value = 28573384
$ws.Range("K62").Value = 2329.5454
$ws.Range("L62").Value = 28573384
$ws.Range("M62").Value = -1705.5454
$ws.Range("N62").Value = -28574632
$ws.Range("H65").Value = 6898791
$ws.Range("I65").Value = 2329.5454
$ws.Range("J65").Value = 28573384
$ws.Range("K65").Value = 11647.727
$ws.Range("L65").Value = 142866920
$ws.Range("M65").Value = -8527.726999999999
$ws.Range("N65").Value = -142873160
$ws.Range("H86").Value = 5176167
$ws.Range("J86").Value = 62801.4
$ws.Range("L86").Value = 62801.4
$ws.Range("N86").Value = -65047.4
$ws.Range("H89").Value = 5176167
$ws.Range("J89").Value = 62801.4
$ws.Range("L89").Value = 314007
$ws.Range("N89").Value = -325239
$ws.Range("H132").Value = 1439.7037
$ws.Range("J132").Value = 2587
$ws.Range("L132").Value = 7761
$ws.Range("N132").Value = -12821
$ws.Range("H134").Value = 25001720
$ws.Range("I134").Value = 1665.8
$ws.Range("K134").Value = 4997.4
$ws.Range("M134").Value = -2462.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 20409116
$ws.Range("J131").Value = 1177.3939
$ws.Range("L131").Value = 3532.1817
$ws.Range("N131").Value = -13612.1817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 26473936
$ws.Range("I70").Value = 20836620
$ws.Range("J70").Value = 40003500
$ws.Range("K70").Value = 20836620
$ws.Range("L70").Value = 40003500
$ws.Range("M70").Value = -20836350
$ws.Range("N70").Value = -40004040
$ws.Range("H73").Value = 26473936
$ws.Range("I73").Value = 20836620
$ws.Range("J73").Value = 40003500
$ws.Range("K73").Value = 20836620
$ws.Range("L73").Value = 40003500
$ws.Range("M73").Value = -20835684
$ws.Range("N73").Value = -40005372
$ws.Range("H108").Value = 59800
$ws.Range("J108").Value = 59800
$ws.Range("L108").Value = 59800
$ws.Range("N108").Value = -67480
$ws.Range("H122").Value = 3390.5264
$ws.Range("I122").Value = 3858.3635
$ws.Range("J122").Value = 2747.25
$ws.Range("K122").Value = 11575.0905
$ws.Range("L122").Value = 8241.75
$ws.Range("M122").Value = -9125.0905
$ws.Range("N122").Value = -13141.75
$ws.Range("H132").Value = 3617.318
$ws.Range("I132").Value = 3635
$ws.Range("J132").Value = 3591.7778
$ws.Range("K132").Value = 10905
$ws.Range("L132").Value = 10775.3334
$ws.Range("M132").Value = -8375
$ws.Range("N132").Value = -15835.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2110.25
$ws.Range("J7").Value = 2399
$ws.Range("L7").Value = 2399
$ws.Range("N7").Value = -2623
$ws.Range("H126").Value = 2110.25
$ws.Range("J126").Value = 2399
$ws.Range("L126").Value = 7197
$ws.Range("N126").Value = -12137
$ws.Range("H132").Value = 2809.2632
$ws.Range("I132").Value = 2355.2144
$ws.Range("K132").Value = 7065.6432
$ws.Range("M132").Value = -4535.6432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2959.3125
$ws.Range("I132").Value = 2654.25
$ws.Range("J132").Value = 3874.5
$ws.Range("K132").Value = 7962.75
$ws.Range("L132").Value = 11623.5
$ws.Range("M132").Value = -5432.75
$ws.Range("N132").Value = -16683.5
$ws.Range("H136").Value = 1278.3448
$ws.Range("I136").Value = 1072.7894
$ws.Range("K136").Value = 3218.3682
$ws.Range("M136").Value = -668.3681999999999

Write-Output "Applied all Kujata_Profits updates"